# Refresh the "cryptos" price table with the latest pull from coinranking.com.
# The source feed re-ranked a few coins (their rows swap places) and every
# Price / Volume(1h) column got new quotes - exactly what's captured below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> { Column letter -> new cell text }.
# Only B (Coin), C (Link), D (Price) and E (Volume(1h)) ever change; A (rank)
# and row 1 (headers) are untouched.
$updates = [ordered]@{
    2  = @{ D = '37.847.93';  E = '  -0.22%  ' }
    3  = @{ D = '2.027.79';   E = '  -1.62%  ' }
    4  = @{               E = '  +0.19%  ' }
    5  = @{ D = '227.02';     E = '  -1.62%  ' }
    6  = @{               E = '  -0.89%  ' }
    7  = @{ D = '59.63';      E = '  +4.09%  ' }
    9  = @{ D = '0.385';      E = '  -0.15%  ' }
    10 = @{ D = '0.0808';     E = '  +0.47%  ' }
    11 = @{               E = '  +0.36%  ' }
    12 = @{ B = 'Chainlink';                    C = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link';                      D = '14.55';      E = '  -0.38%  ' }
    13 = @{ B = 'WrappedliquidstakedEther2.0';  C = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth';       D = '2.328.14';   E = '  -1.53%  ' }
    14 = @{ D = '21.00';      E = '  +1.75%  ' }
    15 = @{ D = '0.751';      E = '  +0.33%  ' }
    16 = @{ D = '5.22';       E = '  -0.89%  ' }
    17 = @{ D = '2.038.06';   E = '  -1.11%  ' }
    18 = @{ D = '37.819.38';  E = '  -0.07%  ' }
    19 = @{ D = '6.03';       E = '  -3.62%  ' }
    20 = @{ D = '69.48';      E = '  -0.34%  ' }
    21 = @{ D = '0.0₃0821';   E = '  -1.12%  ' }
    22 = @{ D = '224.59';     E = '  -0.05%  ' }
    23 = @{               E = '  +0.13%  ' }
    24 = @{ D = '2.42';       E = '  -1.53%  ' }
    25 = @{ D = '2.21';       E = '  -2.56%  ' }
    26 = @{ D = '165.08';     E = '  -0.45%  ' }
    27 = @{ D = '9.17';       E = '  -1.21%  ' }
    28 = @{               E = '  -3.69%  ' }
    29 = @{ D = '18.87';      E = '  -1.36%  ' }
    30 = @{ D = '1.28';       E = '  -6.03%  ' }
    31 = @{ D = '0.120';      E = '  +1.17%  ' }
    32 = @{ D = '4.43';       E = '  -2.50%  ' }
    33 = @{               E = '  +0.77%  ' }
    34 = @{ D = '0.0601';     E = '  -2.17%  ' }
    35 = @{ D = '4.48';       E = '  -1.76%  ' }
    36 = @{               E = '  +6.48%  ' }
    37 = @{ D = '2.24';       E = '  -5.66%  ' }
    38 = @{ D = '3.25';       E = '  -1.65%  ' }
    39 = @{               E = '  -0.08%  ' }
    40 = @{ D = '1.539.85';   E = '  +3.77%  ' }
    41 = @{               E = '  -0.75%  ' }
    42 = @{ D = '96.35';      E = '  -2.63%  ' }
    43 = @{ D = '16.57';      E = '  -1.12%  ' }
    44 = @{ D = '2.82';       E = '  -1.00%  ' }
    45 = @{ D = '0.0919';     E = '  -3.79%  ' }
    46 = @{               E = '  -2.08%  ' }
    47 = @{ B = 'FTXToken';  C = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt';           D = '3.90';  E = '  -5.11%  ' }
    48 = @{ B = 'MXToken';   C = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx';          D = '2.96';  E = '  +0.59%  ' }
    49 = @{ B = 'ARBITRUM';  C = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb';           D = '1.00';  E = '  -1.95%  ' }
    50 = @{ D = '7.07';       E = '  -1.07%  ' }
    51 = @{ D = '2.215.79';   E = '  -1.32%  ' }
}

foreach ($rowNum in $updates.Keys) {
    $rowChanges = $updates[$rowNum]
    foreach ($col in $rowChanges.Keys) {
        $text = $rowChanges[$col]
        $cell = $ws.Range("$col$rowNum")

        # The Price column holds plain-decimal-looking text ("227.02", "1.00", ...).
        # These cells carry no explicit text number format, so a bare .Value write
        # would let Excel reinterpret them as numbers and drop the trailing zeros.
        # Prefix with an apostrophe (exactly what typing it in by hand does) whenever
        # the new text would otherwise be parsed as a plain number (single decimal point).
        $looksNumeric = ($col -eq 'D') -and ($text -match '^[0-9]+(\.[0-9]+)?$')

        if ($looksNumeric) {
            $cell.Value = "'" + $text
        } else {
            $cell.Value = $text
        }
    }
}
